# Actualizar fixtures de todas las ligas (goles por tiempo, stats) 2025-08-15
#
# 1) Add 4 new tracking columns (T:W) with headers, styled like the existing
#    header row (bold / centered / bordered, same as S1).
# 2) Append 8 new fixture rows (154-161) with full match data; the two new
#    "posesion" columns (T,U) are left blank for these rows (not yet
#    available), while V/W record the ingestion source/status.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) New header cells T1:U1:V1:W1 -----------------------------------
# Copy S1's style (bold, centered, thin border) onto the new header cells
# first, then set their text so the style carries over cleanly.
$ws.Range("S1").Copy($ws.Range("T1:W1"))
$ws.Range("T1").Value = "Posesión Local (%)"
$ws.Range("U1").Value = "Posesión Visita (%)"
$ws.Range("V1").Value = "fuente_tiempos"
$ws.Range("W1").Value = "estado_datos"

# --- 2) New fixture rows 154-161 ---------------------------------------
# Column A holds dates as plain text in this sheet (not date serials), so
# force text format on the new rows before writing to avoid Excel's
# automatic date-string -> serial-number conversion.
$ws.Range("A154:A161").NumberFormat = "@"

# Fecha|Local|Visita|Goles Local|Goles Visita|Fixture ID|Corners Local|Corners Visita|
# Amarillas Local|Amarillas Visita|Rojas Local|Rojas Visita|Goles 1T Local|Goles 1T Visita|
# Goles 2T Local|Goles 2T Visita|Posesion Local|Posesion Visita|Resultado
$newRows = @"
2025-08-08|San Luis|Magallanes|0|1|1348387|9|5|2|3|1|1|0|0|0|1|50|50|V
2025-08-09|Rangers de Talca|Santiago Morning|0|0|1348391|1|11|5|4|0|0|0|0|0|0|48|52|E
2025-08-09|Deportes Santa Cruz|Universidad de Concepcion|1|2|1348390|7|6|2|5|0|1|1|2|0|0|55|45|V
2025-08-09|Santiago Wanderers|Deportes Copiapo|1|0|1348388|4|3|1|3|0|3|0|0|1|0|51|49|L
2025-08-10|San Marcos de Arica|Deportes Temuco|0|0|1348385|6|5|4|4|0|0|0|0|0|0|64|36|E
2025-08-10|Antofagasta|Cobreloa|3|0|1348386|2|7|4|1|0|1|2|0|1|0|46|54|L
2025-08-10|Recoleta|Union San Felipe|2|1|1348389|4|10|3|2|0|0|1|0|1|1|45|55|L
2025-08-10|Concepción|Curico Unido|0|0|1348392|7|2|4|5|0|2|0|0|0|0|55|45|E
"@

$lines = $newRows.Trim() -split "`n"
$r = 154
foreach ($line in $lines) {
    $p = $line.Trim() -split '\|'

    $ws.Cells.Item($r, 1).Value = $p[0]          # A Fecha
    $ws.Cells.Item($r, 2).Value = $p[1]          # B Local
    $ws.Cells.Item($r, 3).Value = $p[2]          # C Visita
    $ws.Cells.Item($r, 4).Value = [double]$p[3]  # D Goles Local
    $ws.Cells.Item($r, 5).Value = [double]$p[4]  # E Goles Visita
    $ws.Cells.Item($r, 6).Value = [double]$p[5]  # F Fixture ID
    $ws.Cells.Item($r, 7).Value = [double]$p[6]  # G Corners Local
    $ws.Cells.Item($r, 8).Value = [double]$p[7]  # H Corners Visita
    $ws.Cells.Item($r, 9).Value = [double]$p[8]  # I Amarillas Local
    $ws.Cells.Item($r, 10).Value = [double]$p[9] # J Amarillas Visita
    $ws.Cells.Item($r, 11).Value = [double]$p[10] # K Rojas Local
    $ws.Cells.Item($r, 12).Value = [double]$p[11] # L Rojas Visita
    $ws.Cells.Item($r, 13).Value = [double]$p[12] # M Goles 1T Local
    $ws.Cells.Item($r, 14).Value = [double]$p[13] # N Goles 1T Visita
    $ws.Cells.Item($r, 15).Value = [double]$p[14] # O Goles 2T Local
    $ws.Cells.Item($r, 16).Value = [double]$p[15] # P Goles 2T Visita
    $ws.Cells.Item($r, 17).Value = [double]$p[16] # Q Posesion Local (%)
    $ws.Cells.Item($r, 18).Value = [double]$p[17] # R Posesion Visita (%)
    $ws.Cells.Item($r, 19).Value = $p[18]         # S Resultado

    # T (Posesión Local %) and U (Posesión Visita %) stay blank for these
    # rows - not yet available - matching the source edit.
    $ws.Cells.Item($r, 22).Value = "score"        # V fuente_tiempos
    $ws.Cells.Item($r, 23).Value = "OK"           # W estado_datos

    $r = $r + 1
}
